$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record values for each data row (2-43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 29).Value = 76   # AC
    $ws.Cells.Item($r, 30).Value = 86   # AD
    $ws.Cells.Item($r, 31).Value = 0    # AE
}

$excel.CutCopyMode = 0
